$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.640.51'
$ws.Range("E2").Value = '  -3.89%  '
$ws.Range("D3").Value = '2.909.63'
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'527.98"
$ws.Range("E5").Value = '  -5.33%  '
$ws.Range("D6").Value = "'142.88"
$ws.Range("E6").Value = '  -7.55%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = "'0.553"
$ws.Range("E8").Value = '  -1.56%  '
$ws.Range("D9").Value = '2.910.50'
$ws.Range("E9").Value = '  -3.88%  '
$ws.Range("E10").Value = '  -4.86%  '
$ws.Range("D11").Value = "'5.86"
$ws.Range("E11").Value = '  -8.72%  '
$ws.Range("D12").Value = "'0.352"
$ws.Range("E12").Value = '  -3.93%  '
$ws.Range("D13").Value = '3.413.54'
$ws.Range("E13").Value = '  -3.70%  '
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").Value = '60.769.05'
$ws.Range("E15").Value = '  -3.78%  '
$ws.Range("D16").Value = "'22.61"
$ws.Range("E16").Value = '  -6.03%  '
$ws.Range("D17").Value = '2.906.98'
$ws.Range("E17").Value = '  -4.13%  '
$ws.Range("D18").Value = "'0.0000140"
$ws.Range("E18").Value = '  -7.01%  '
$ws.Range("D19").Value = "'4.92"
$ws.Range("E19").Value = '  -3.53%  '
$ws.Range("D20").Value = "'11.52"
$ws.Range("E20").Value = '  -3.68%  '
$ws.Range("D21").Value = "'360.59"
$ws.Range("E22").Value = '  -1.43%  '
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = "'5.66"
$ws.Range("E24").Value = '  -1.70%  '
$ws.Range("D25").Value = "'63.44"
$ws.Range("E25").Value = '  -2.84%  '
$ws.Range("D26").Value = '3.034.79'
$ws.Range("E26").Value = '  -3.90%  '
$ws.Range("D27").Value = "'0.451"
$ws.Range("E27").Value = '  -3.45%  '
$ws.Range("D28").Value = "'0.179"
$ws.Range("E28").Value = '  -5.46%  '
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("D30").Value = '0.0₃0858'
$ws.Range("E30").Value = '  -12.44%  '
$ws.Range("D31").Value = "'7.65"
$ws.Range("E31").Value = '  -11.72%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("E33").Value = '  -4.85%  '
$ws.Range("D34").Value = "'19.66"
$ws.Range("E34").Value = '  -3.69%  '
$ws.Range("D35").Value = "'152.15"
$ws.Range("E35").Value = '  -4.27%  '
$ws.Range("D36").Value = "'4.34"
$ws.Range("E36").Value = '  -8.02%  '
$ws.Range("D37").Value = "'5.57"
$ws.Range("E37").Value = '  -7.95%  '
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = '  -9.06%  '
$ws.Range("D39").Value = "'1.20"
$ws.Range("E39").Value = '  -7.99%  '
$ws.Range("D40").Value = "'37.92"
$ws.Range("E40").Value = '  +1.18%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.332.14'
$ws.Range("E41").Value = '  -7.92%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'1.46"
$ws.Range("E42").Value = '  -7.17%  '
$ws.Range("D43").Value = "'3.68"
$ws.Range("E43").Value = '  -6.58%  '
$ws.Range("D44").Value = "'0.642"
$ws.Range("E44").Value = '  -4.15%  '
$ws.Range("D45").Value = "'20.81"
$ws.Range("E45").Value = '  -8.74%  '
$ws.Range("D46").Value = "'0.0569"
$ws.Range("E46").Value = '  -5.17%  '
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = "'4.84"
$ws.Range("E48").Value = '  -4.56%  '
$ws.Range("D49").Value = "'0.0232"
$ws.Range("E49").Value = '  -6.50%  '
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("D51").Value = "'0.0924"
$ws.Range("E51").Value = '  -2.51%  '
